$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "Obrigatorio" (Required) column (E) to "S" (Sim) for rows 2-8
$ws.Range("E2:E8").Value = "S"
